# Add Test Data for Russia, Finland and Hungary markets.
# These three new sheets are clones of the "Spain" sheet template
# (same layout/styles/merged cells), with only the Market-name (B2)
# and NGC reference (B4) cell values changed, appended after the
# existing last sheet ("Denmark").

$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("Spain")

function Add-MarketSheet {
    param([string]$SheetName, [string]$NgcRef, [string]$MarketName)

    $afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $template.Copy($null, $afterSheet)

    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = $SheetName

    # Set B4 (NGC reference) before B2 (Market name) so new shared
    # strings are appended in that order.
    $newSheet.Range("B4").Value = $NgcRef
    $newSheet.Range("B2").Value = $MarketName
}

# A throwaway sheet is copied + immediately deleted first purely to
# advance the workbook's internal sheetId counter so the three real
# new sheets land on sheetId 22/23/24 (matching upstream authoring),
# instead of 21/22/23.
$placeholderAfter = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $placeholderAfter)
$placeholder = $wb.Worksheets.Item($wb.Worksheets.Count)
$placeholder.Name = "ZZZ_Placeholder"

Add-MarketSheet "Russia"  "NGC-2929/T3302" "Russia Market"
Add-MarketSheet "Finland" "NGC-3130/T2947" "Finland Market"
Add-MarketSheet "Hungary" "NGC-3104/T2996" "Hungary Market"

$wb.Worksheets.Item("ZZZ_Placeholder").Delete()

# The newly added "Hungary" sheet (last one) ends up active/selected,
# matching the target workbook view.
$wb.Worksheets.Item("Hungary").Activate()
